$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(15, 1).Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Cells.Item(15, 2).Value = "all"
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = "2025-08-06 23:45:52"
$ws.Cells.Item(15, 5).Value = "Paid"

$ws.Cells.Item(16, 1).Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Cells.Item(16, 2).Value = "all"
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = "2025-08-06 23:59:57"
$ws.Cells.Item(16, 5).Value = "Paid"
